# The first data row ("H 72") was removed from the table; every row
# below it shifts up by one (last row 63 disappears, new last row is 62).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
